$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout (before):
#   A1=from   B1=to
#   A2=0.591...  B2=0.580...
#
# Target layout (after):
#   A1=Area       B1=from      C1=to
#   A2=California B2=0.591...  C2=0.580...
#   A3=New York   B3=0.690...  C3=0.680...

# Insert a new column before column A to make room for the "Area" column.
$ws.Range("A:A").Insert()

# Fill in the header row.
$ws.Range("A1").Value = "Area"

# Fill in the existing row's new "Area" value.
$ws.Range("A2").Value = "California"

# Add the new row for New York.
$ws.Range("A3").Value = "New York"
$ws.Range("B3").Value = 0.69032135242221404
$ws.Range("C3").Value = 0.68075655255527501

# Update selection to match target state.
$ws.Range("I20").Select()
